# "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had three stray "header/footer only" rows interleaved with the
# real data (they carry a label in column A but no B:H values):
#   row 5  -> "situação do domicílio"
#   row 8  -> "grandes regiões e unidades da federação"
#   row 41 -> "fonte: ibge, diretoria de pesquisas, ..." (footer note)
#
# Removing those rows re-aligns every region label with its correct data
# row (the data itself is untouched - it just collapses upward by one slot
# at each removed separator). We also fix the stray pandas leftover header
# text in B2 ("unnamed: 1_level_1" -> "total").
#
# Rows are deleted bottom-to-top so that each row number used below still
# refers to the same original row (deleting a row only shifts the rows
# that come after it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").EntireRow.Delete()
$ws.Range("A8").EntireRow.Delete()
$ws.Range("A5").EntireRow.Delete()

$ws.Range("B2").Value = "total"
